$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.132.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.748.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'236.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.5291"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("D8").Value = "'0.2833"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.06182"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "1.746.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'0.07175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "'15.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'0.6460"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'4.629"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").Value = "'78.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'0.9996"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "26.025.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "'11.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").Value = "'0.000006744"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").Value = "1.968.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'4.329"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("D23").Value = "'8.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'5.230"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "'139.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").Value = "'15.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'1.805"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "'104.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "'0.08310"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'3.798"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("D32").Value = "'3.640"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.19%  "
$ws.Range("D33").Value = "'0.04634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").Value = "'2.645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").Value = "'0.6346"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("D37").Value = "'2.703"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "'0.01626"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("D39").Value = "'1.974"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("D40").Value = "'0.9994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "'102.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").Value = "'0.3929"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'0.7499"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "'5.057"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "'0.1153"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").Value = "'6.357"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.05344"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").Value = "'54.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.58%  "
$ws.Range("D49").Value = "'31.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.74%  "
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").Value = "'7.624"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.37%  "
